$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Establish the "data row" look (Calibri 10pt, theme Text-1 colour, kept
#     as Text so numeric-looking values like "1985" stay strings) on A2,
#     then fan that formatting out to the rest of the two new rows. ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").Font.Size = 10
$ws.Range("A2").Font.ThemeColor = 1
$ws.Range("A2").Value = "MCH169-1"

$ws.Range("A2").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)
$ws.Range("C2:H3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 2 (left to right) ---
$ws.Range("C2").Value = "DIARIES 1985 I, 1985 II, 1986 I"
$ws.Range("D2").Value = "1985"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 22C | GRAP COUNT NUMER: NONE"

# --- Row 3 (left to right) ---
$ws.Range("A3").Value = "MCH169-2"
$ws.Range("C3").Value = "1986 II, 1987 I, 1987 II"
$ws.Range("D3").Value = "1986"
$ws.Range("E3").Value = "Series"
$ws.Range("F3").Value = "1 Box"
$ws.Range("G3").Value = "LOCATION: 22C | GRAP COUNT NUMER: NONE"

# --- Restore the frozen header pane / selection over the new data extent ---
$ws.Range("A2:I3").Select()
$excel.ActiveWindow.FreezePanes = $true
